$d = $word.ActiveDocument

# Build the shared run-properties / paragraph-mark fragment used by the
# new paragraphs (Times New Roman, 14pt / half-points 28).
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Paragraph 1: "Путь 1: " + "1-2-3-4-5-6-7-8-9-15" kept as two separate
# runs (identical formatting, but authored as two distinct insertions).
$para1 = '<w:p><w:pPr>' + $rPr + '</w:pPr>' `
       + '<w:r>' + $rPr + '<w:t xml:space="preserve">Путь 1: </w:t></w:r>' `
       + '<w:r>' + $rPr + '<w:t>1-2-3-4-5-6-7-8-9-15</w:t></w:r>' `
       + '</w:p>'

# Paragraph 2: "Путь 2: 1-2-3-4-5-6-7-8-9-10-11-12-13-14" as a single run.
$para2 = '<w:p><w:pPr>' + $rPr + '</w:pPr>' `
       + '<w:r>' + $rPr + '<w:t>Путь 2: 1-2-3-4-5-6-7-8-9-10-11-12-13-14</w:t></w:r>' `
       + '</w:p>'

# Insert both new paragraphs at the very end of the document (after the
# last existing paragraph, before the final section break).
$r1 = $d.Content
$r1.Collapse(0)
$r1.InsertXML($pkgHeader + $para1 + $pkgFooter)

$r2 = $d.Content
$r2.Collapse(0)
$r2.InsertXML($pkgHeader + $para2 + $pkgFooter)
